# Applies the edit described by the diff:
#  - Fix JOB_TITLE typo in A10: "System Analyst 1234" -> "System Analyst"
#  - Fix NATIONALITY in J10: "51 Kuwaiti" -> "Syria"
#  - Move the active cell selection to J11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the job title that had a stray "1234" suffix.
$ws.Range("A10").Value = "System Analyst"

# Correct the nationality value for this row.
$ws.Range("J10").Value = "Syria"

# Update the current selection/active cell as recorded in the saved view state.
$ws.Range("J11").Select()
